{"js": "// Replace the date line and every \"A\u00d7B=\" multiplication prompt in the\n// table with the new values from the target revision. Every old string\n// below is unique within the document, so a scoped search+replace per\n// pair is unambiguous.\nconst replacements = [\n  [\"2024-07-06 Saturday\", \"2024-07-07 Sunday\"],\n  [\"56\u00d765=\", \"39\u00d735=\"],\n  [\"84\u00d780=\", \"48\u00d756=\"],\n  [\"58\u00d740=\", \"60\u00d735=\"],\n  [\"19\u00d765=\", \"22\u00d759=\"],\n  [\"45\u00d762=\", \"97\u00d721=\"],\n  [\"49\u00d714=\", \"89\u00d781=\"],\n  [\"15\u00d770=\", \"94\u00d788=\"],\n  [\"95\u00d711=\", \"60\u00d748=\"],\n  [\"70\u00d799=\", \"76\u00d721=\"],\n  [\"70\u00d792=\", \"76\u00d722=\"],\n  [\"98\u00d771=\", \"88\u00d737=\"],\n  [\"91\u00d798=\", \"41\u00d794=\"],\n  [\"46\u00d741=\", \"11\u00d734=\"],\n  [\"38\u00d713=\", \"49\u00d761=\"],\n  [\"66\u00d793=\", \"77\u00d779=\"],\n  [\"73\u00d736=\", \"17\u00d713=\"],\n  [\"85\u00d741=\", \"36\u00d740=\"],\n  [\"55\u00d726=\", \"87\u00d766=\"],\n  [\"90\u00d767=\", \"89\u00d734=\"],\n  [\"30\u00d760=\", \"53\u00d718=\"],\n  [\"68\u00d740=\", \"55\u00d715=\"],\n  [\"51\u00d787=\", \"28\u00d785=\"],\n  [\"26\u00d792=\", \"45\u00d742=\"],\n  [\"68\u00d750=\", \"12\u00d719=\"],\n  [\"13\u00d713=\", \"78\u00d769=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n\n", "ps1": "# Replace the date line and every \"A\u00d7B=\" multiplication prompt in the\n# table with the new values from the target revision, using Word's\n# Find/Replace (wdReplaceAll) over the whole document. Every \"old\" string\n# is unique in the document, so each call only ever touches the one cell\n# it targets.\n$d = $word.ActiveDocument\n\nfunction Replace-Text($old, $new) {\n  $find = $d.Content.Find\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n\nReplace-Text \"2024-07-06 Saturday\" \"2024-07-07 Sunday\"\nReplace-Text \"56\u00d765=\" \"39\u00d735=\"\nReplace-Text \"84\u00d780=\" \"48\u00d756=\"\nReplace-Text \"58\u00d740=\" \"60\u00d735=\"\nReplace-Text \"19\u00d765=\" \"22\u00d759=\"\nReplace-Text \"45\u00d762=\" \"97\u00d721=\"\nReplace-Text \"49\u00d714=\" \"89\u00d781=\"\nReplace-Text \"15\u00d770=\" \"94\u00d788=\"\nReplace-Text \"95\u00d711=\" \"60\u00d748=\"\nReplace-Text \"70\u00d799=\" \"76\u00d721=\"\nReplace-Text \"70\u00d792=\" \"76\u00d722=\"\nReplace-Text \"98\u00d771=\" \"88\u00d737=\"\nReplace-Text \"91\u00d798=\" \"41\u00d794=\"\nReplace-Text \"46\u00d741=\" \"11\u00d734=\"\nReplace-Text \"38\u00d713=\" \"49\u00d761=\"\nReplace-Text \"66\u00d793=\" \"77\u00d779=\"\nReplace-Text \"73\u00d736=\" \"17\u00d713=\"\nReplace-Text \"85\u00d741=\" \"36\u00d740=\"\nReplace-Text \"55\u00d726=\" \"87\u00d766=\"\nReplace-Text \"90\u00d767=\" \"89\u00d734=\"\nReplace-Text \"30\u00d760=\" \"53\u00d718=\"\nReplace-Text \"68\u00d740=\" \"55\u00d715=\"\nReplace-Text \"51\u00d787=\" \"28\u00d785=\"\nReplace-Text \"26\u00d792=\" \"45\u00d742=\"\nReplace-Text \"68\u00d750=\" \"12\u00d719=\"\nReplace-Text \"13\u00d713=\" \"78\u00d769=\"\n"}
